# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text in A1 ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $hoja1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 1.85 = 6854.04 pesos"), "1000 Bs = 1.85 = 6856.02 pesos"
$text = $text -replace [regex]::Escape("6854.04 pesos = 1.84 = 879.24 Bs"), "6856.02 pesos = 1.83 = 887.72 Bs"
$cell.Value2 = $text

# --- tasas: update the numeric rate cells ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("O10").Value = 3707.05
$tasas.Range("N12").Value = 3738
$tasas.Range("O12").Value = 484
